# Auto-generated Excel COM-interop script applying the diff changes
# to "Jogos_da_Semana_FlashScore_2024-10-16.xlsx" (odds update for 2024-10-16).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 4.5
$ws.Range("I4").Value = 1.85
$ws.Range("L4").Value = 2.6
$ws.Range("AE4").Value = 21
$ws.Range("AF4").Value = 81
$ws.Range("AI4").Value = 7.5
$ws.Range("AJ4").Value = 9
$ws.Range("AU4").Value = 9.5
$ws.Range("G5").Value = 1.57
$ws.Range("H5").Value = 3.9
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 10
$ws.Range("O5").Value = 1.36
$ws.Range("P5").Value = 3
$ws.Range("Q5").Value = 2.15
$ws.Range("R5").Value = 1.67
$ws.Range("AU5").Value = 9.5
$ws.Range("H6").Value = 3.8
$ws.Range("J6").Value = 6
$ws.Range("K6").Value = 2.2
$ws.Range("Q6").Value = 2.05
$ws.Range("R6").Value = 1.75
$ws.Range("U6").Value = 2.1
$ws.Range("V6").Value = 1.67
$ws.Range("W6").Value = 15
$ws.Range("AC6").Value = 8.5
$ws.Range("AD6").Value = 7
$ws.Range("AE6").Value = 19
$ws.Range("AF6").Value = 67
$ws.Range("AH6").Value = 6
$ws.Range("AJ6").Value = 8.5
$ws.Range("AQ6").Value = 126
$ws.Range("AU6").Value = 9
$ws.Range("AY6").Value = 21
$ws.Range("BB6").Value = 151
$ws.Range("N8").Value = 8
$ws.Range("O8").Value = 1.31
$ws.Range("P8").Value = 2.9
$ws.Range("Q8").Value = 1.95
$ws.Range("R8").Value = 1.75
$ws.Range("U8").Value = 1.7
$ws.Range("V8").Value = 1.91
$ws.Range("W8").Value = 9
$ws.Range("Y8").Value = 11
$ws.Range("AA8").Value = 28
$ws.Range("AC8").Value = 9.25
$ws.Range("AF8").Value = 65
$ws.Range("AG8").Value = 500
$ws.Range("AH8").Value = 7.7
$ws.Range("AI8").Value = 11
$ws.Range("AL8").Value = 17.5
$ws.Range("AM8").Value = 27
$ws.Range("AU8").Value = 6.7
$ws.Range("BA8").Value = 65
$ws.Range("M9").Value = 1.07
$ws.Range("N9").Value = 9
$ws.Range("G12").Value = 1.95
$ws.Range("H12").Value = 3.2
$ws.Range("I12").Value = 3.9
$ws.Range("J12").Value = 2.63
$ws.Range("L12").Value = 4.33
$ws.Range("M12").Value = 1.07
$ws.Range("N12").Value = 9
$ws.Range("R12").Value = 1.8
$ws.Range("S12").Value = 1.4
$ws.Range("T12").Value = 2.75
$ws.Range("AD12").Value = 6.5
$ws.Range("AT12").Value = 2.75

Write-Host "Applied 68 cell updates."
